$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "wafer"
$ws.Range("A2").Value = "P01"
[void]$ws.Range("D2").Select()
